$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph. We rebuild paragraph 1's whole range
#    (title + a brand-new Normal paragraph) via InsertXML so the new
#    paragraph comes out with the exact run/style layout Word itself
#    would produce - no stray pPr, no leftover rsid noise.
# ------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$titlePara = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Play Dragon''s Inferno for Free - Review and Gameplay Details</w:t></w:r></w:p>'
$metaPara  = '<w:p ' + $wNs + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dragon''s Inferno, a popular slot game. Play for free and learn about the game''s features, symbols, and more.</w:t></w:r></w:p>'

$p1 = $d.Paragraphs.Item(1)
[void]$p1.Range.InsertXML($titlePara + $metaPara)

# ------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Dragon's Inferno..." paragraph
#    that used to sit right before the closing italic paragraph.
#    (Paragraph.Range.Text includes the trailing paragraph mark, so
#    trim it before comparing; and skip the real Heading1 title.)
# ------------------------------------------------------------------
$old = "Play Dragon's Inferno for Free - Review and Gameplay Details"
For ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13, [char]7)
    If ($paraText -eq $old -and $para.Style.NameLocal -ne "Heading 1") {
        $para.Range.Delete()
    }
}

# ------------------------------------------------------------------
# 3. Swap the closing italic paragraph's text for the new image-prompt
#    copy, keeping its existing (italic) run formatting intact. We set
#    Range.Text directly (rather than Find/Replace) so Word's smart-
#    quote autocorrect doesn't curl the apostrophes.
# ------------------------------------------------------------------
$newClosing = "Create a feature image for Dragon's Inferno that features a happy Maya warrior wearing glasses in a cartoon style. The warrior should be holding a dragon's egg in one hand and a treasure chest overflowing with gold coins in the other. In the background, there should be a fiery landscape with dragons flying in the distance. The image should be bright and colorful, capturing the adventurous spirit and fun gameplay of this slot machine game."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$closingRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$closingRange.Text = $newClosing

Write-Output "done"
